# Added 1.1.0 of term
# - bump the "Version" metadata value 1.0.0 -> 1.1.0
# - bump the "Date" metadata value to the new publish timestamp
# - re-assert the existing top/wrap alignment on the formatted ranges so the
#   cellXfs records explicitly carry applyAlignment="true" (matches the
#   upstream re-save that flips applyAlignment on for the header/body styles)

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "1.1.0"
$wsMeta.Range("B8").Value = "2023-07-10T23:08:03+02:00"

$wsInclude = $wb.Worksheets.Item("Include from FSIII")

# Only touch the cells that actually already hold data (mirrors the
# upstream re-save, which flips applyAlignment on for the two existing
# styles without materialising any new, previously-empty cells).
$wsMeta.Range("A1:B14").VerticalAlignment = -4160  # xlTop
$wsMeta.Range("A1:B14").WrapText = $true

$wsInclude.Range("A1:C2").VerticalAlignment = -4160  # xlTop
$wsInclude.Range("A1:C2").WrapText = $true
$wsInclude.Range("A3:B4").VerticalAlignment = -4160  # xlTop
$wsInclude.Range("A3:B4").WrapText = $true
